$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new cabinet "UATS" below the last existing row (A14 -> A15)
$ws.Range("A15").Value = "UATS"

# Copy the formatting of the previous data row so the new entry matches
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122) | Out-Null

# Update selection to mirror the post-entry cursor position (A25)
$ws.Range("A25").Select() | Out-Null
